$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 2 de Julio de 2020 a las 18:09"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 2799226
$ws.Range("C4").Value = 19273
$ws.Range("D4").Value = 1168436
$ws.Range("E4").Value = 1499806
$ws.Range("G4").Value = 186
$ws.Range("H4").Value = 130984

# Row 7: India -> India
$ws.Range("B7").Value = 612486
$ws.Range("C7").Value = 7266
$ws.Range("D7").Value = 366027
$ws.Range("E7").Value = 228463
$ws.Range("G7").Value = 148
$ws.Range("H7").Value = 17996

# Row 11: Chile -> Chile
$ws.Range("B11").Value = 284541
$ws.Range("C11").Value = 2498
$ws.Range("D11").Value = 249247
$ws.Range("E11").Value = 29374
$ws.Range("G11").Value = 167
$ws.Range("H11").Value = 5920

# Row 12: Italia -> Italia
$ws.Range("B12").Value = 240961
$ws.Range("C12").Value = 201
$ws.Range("D12").Value = 191083
$ws.Range("E12").Value = 15060
$ws.Range("G12").Value = 30
$ws.Range("H12").Value = 34818

# Row 38: Singapur -> Singapur
$ws.Range("D38").Value = 39429
$ws.Range("E38").Value = 4855

# Row 46: Republica Dominicana -> Republica Dominicana
$ws.Range("B46").Value = 34197
$ws.Range("C46").Value = 810
$ws.Range("D46").Value = 18141
$ws.Range("E46").Value = 15291
$ws.Range("G46").Value = 11
$ws.Range("H46").Value = 765

# Row 58: Ghana -> Azerbaiyan
$ws.Range("A58").Value = "Azerbaiyan"
$ws.Range("B58").Value = 18684
$ws.Range("C58").Value = 572
$ws.Range("D58").Value = 10425
$ws.Range("E58").Value = 8031
$ws.Range("G58").Value = 8
$ws.Range("H58").Value = 228

# Row 59: Azerbaiyan -> Ghana
$ws.Range("A59").Value = "Ghana"
$ws.Range("B59").Value = 18134
$ws.Range("D59").Value = 13550
$ws.Range("E59").Value = 4467
$ws.Range("H59").Value = 117

# Row 61: Moldavia -> Moldavia
$ws.Range("B61").Value = 17150
$ws.Range("C61").Value = 252
$ws.Range("E61").Value = 6744
$ws.Range("G61").Value = 11
$ws.Range("H61").Value = 560

# Row 63: Nepal -> Argelia
$ws.Range("A63").Value = "Argelia"
$ws.Range("B63").Value = 14657
$ws.Range("C63").Value = 385
$ws.Range("D63").Value = 10040
$ws.Range("E63").Value = 3689
$ws.Range("G63").Value = 8
$ws.Range("H63").Value = 928

# Row 64: Argelia -> Nepal
$ws.Range("A64").Value = "Nepal"
$ws.Range("B64").Value = 14519
$ws.Range("C64").Value = 473
$ws.Range("D64").Value = 5320
$ws.Range("E64").Value = 9168
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 31

# Row 69: Chequia -> Chequia
$ws.Range("B69").Value = 12116
$ws.Range("C69").Value = 70
$ws.Range("D69").Value = 7821
$ws.Range("E69").Value = 3942
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = 353

# Row 80: Kenia -> Kenia
$ws.Range("D80").Value = 2109
$ws.Range("E80").Value = 4680
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = 152

# Row 91: Republica de Yibuti -> Republica de Yibuti
$ws.Range("B91").Value = 4715
$ws.Range("C91").Value = 11
$ws.Range("D91").Value = 4564
$ws.Range("E91").Value = 96

# Row 94: Luxemburgo -> Luxemburgo
$ws.Range("B94").Value = 4395
$ws.Range("C94").Value = 50
$ws.Range("D94").Value = 4012
$ws.Range("E94").Value = 273

# Row 98: Grecia -> Grecia
$ws.Range("B98").Value = 3458
$ws.Range("C98").Value = 26
$ws.Range("E98").Value = 1892

# Row 111: Sri Lanka -> Sri Lanka
$ws.Range("B111").Value = 2066
$ws.Range("C111").Value = 12
$ws.Range("E111").Value = 228

# Row 130: Tunez -> Tunez
$ws.Range("B130").Value = 1178
$ws.Range("C130").Value = 3
$ws.Range("D130").Value = 1039
$ws.Range("E130").Value = 89

# Row 131: Jordania -> Jordania
$ws.Range("B131").Value = 1136
$ws.Range("C131").Value = 3
$ws.Range("D131").Value = 889

# Row 145: Liberia -> Liberia
$ws.Range("B145").Value = 819
$ws.Range("C145").Value = 15
$ws.Range("D145").Value = 338
$ws.Range("E145").Value = 444

# Row 152: Zimbabue -> Montenegro
$ws.Range("A152").Value = "Montenegro"
$ws.Range("B152").Value = 616
$ws.Range("C152").Value = 40
$ws.Range("D152").Value = 315
$ws.Range("E152").Value = 289
$ws.Range("H152").Value = 12

# Row 153: Montenegro -> Zimbabue
$ws.Range("A153").Value = "Zimbabue"
$ws.Range("B153").Value = 605
$ws.Range("D153").Value = 166
$ws.Range("E153").Value = 432
$ws.Range("H153").Value = 7

# Row 205: Fiyi -> Dominica
$ws.Range("A205").Value = "Dominica"

# Row 206: Dominica -> Fiyi
$ws.Range("A206").Value = "Fiyi"
